{"js": "// Convert the Word field `{ m:'doc.html'.fromHTMLURI() }` (stored as\n// fldChar/instrText runs) into plain literal text runs that spell out\n// the same content using curly braces: \"{m:'doc.html'.fromHTMLURI()}\".\n// The _GoBack bookmark that sits in the middle of the field code is kept\n// exactly where it was, between \"doc.html\" and the closing quote.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Locate the paragraph that holds the field (begin/separate/end fldChar\n// runs) by inspecting each paragraph's OOXML.\nconst ooxmlResults = paragraphs.items.map((p) => p.getOoxml());\nawait context.sync();\n\nlet targetParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const xml = ooxmlResults[i].value;\n  if (xml && xml.indexOf(\"fldChar\") !== -1) {\n    targetParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (targetParagraph) {\n  // Build the replacement runs (plain text, no field characters) as a\n  // standalone paragraph fragment and insert it at the very start of the\n  // paragraph that contains the field. Inserting at \"Start\" (rather than\n  // \"Replace\") keeps the existing paragraph mark/properties untouched.\n  const replacementOoxml =\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' +\n    '<w:p>' +\n    '<w:r><w:t>{</w:t></w:r>' +\n    '<w:r><w:t>m</w:t></w:r>' +\n    '<w:r><w:t>:</w:t></w:r>' +\n    \"<w:r><w:t>'</w:t></w:r>\" +\n    '<w:r><w:t>doc.html</w:t></w:r>' +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n    '<w:bookmarkEnd w:id=\"0\"/>' +\n    \"<w:r><w:t>'.fromHTMLURI()</w:t></w:r>\" +\n    '<w:r><w:t xml:space=\"preserve\">}</w:t></w:r>' +\n    '</w:p>' +\n    '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>';\n\n  targetParagraph.insertOoxml(replacementOoxml, \"Start\");\n  await context.sync();\n\n  // Now remove the original field (fldChar begin/separate/end and the\n  // instrText runs). The field's own bookmark gets removed with it, but\n  // the _GoBack bookmark we just inserted above stays in place.\n  const fields = body.fields;\n  fields.load(\"items\");\n  await context.sync();\n\n  if (fields.items.length > 0) {\n    fields.items[0].delete();\n    await context.sync();\n  }\n}\n", "ps1": "# Convert the Word field `{ m:'doc.html'.fromHTMLURI() }` (stored as\n# fldChar/instrText runs) into plain literal text runs that spell out\n# the same content using curly braces: \"{m:'doc.html'.fromHTMLURI()}\".\n# The _GoBack bookmark that sits in the middle of the field code is kept\n# exactly where it was, between \"doc.html\" and the closing quote.\n\n$d = $word.ActiveDocument\n\nif ($d.Fields.Count -gt 0) {\n    $f = $d.Fields.Item(1)\n    $codeStart = $f.Code.Start\n\n    # Locate the paragraph that contains the field.\n    $targetParagraph = $null\n    foreach ($p in $d.Paragraphs) {\n        $s = $p.Range.Start\n        $e = $p.Range.End\n        if ($codeStart -ge $s -and $codeStart -lt $e) {\n            $targetParagraph = $p\n            break\n        }\n    }\n\n    if ($targetParagraph -ne $null) {\n        $quote = [char]39\n\n        # Build the replacement runs (plain text, no field characters) as a\n        # standalone paragraph fragment and insert it at the very start of\n        # the paragraph that contains the field. Inserting a collapsed\n        # range at the paragraph start (rather than replacing the whole\n        # paragraph range) keeps the existing paragraph mark/properties\n        # untouched.\n        $xml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n            '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n            '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n            '<pkg:xmlData>' +\n            '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n            '<w:body>' +\n            '<w:p>' +\n            '<w:r><w:t>{</w:t></w:r>' +\n            '<w:r><w:t>m</w:t></w:r>' +\n            '<w:r><w:t>:</w:t></w:r>' +\n            '<w:r><w:t>' + $quote + '</w:t></w:r>' +\n            '<w:r><w:t>doc.html</w:t></w:r>' +\n            '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n            '<w:bookmarkEnd w:id=\"0\"/>' +\n            '<w:r><w:t>' + $quote + '.fromHTMLURI()</w:t></w:r>' +\n            '<w:r><w:t xml:space=\"preserve\">}</w:t></w:r>' +\n            '</w:p>' +\n            '</w:body>' +\n            '</w:document>' +\n            '</pkg:xmlData>' +\n            '</pkg:part>' +\n            '</pkg:package>'\n\n        $insertPoint = $d.Range($targetParagraph.Range.Start, $targetParagraph.Range.Start)\n        $insertPoint.InsertXML($xml)\n\n        # Now remove the original field (fldChar begin/separate/end and the\n        # instrText runs). The field's own bookmark gets removed with it,\n        # but the _GoBack bookmark we just inserted above stays in place.\n        $d.Fields.Item(1).Delete()\n    }\n}\n"}
